# Swap the 5 doctor names shown on the sheet for the first few entries
# in the shared-strings table (prepping fixture data for the grouping
# feature work).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Dr. K.A. Mohan"
$ws.Range("A2").Value = "Dr. Achuth M Baliga"
$ws.Range("A3").Value = "Dr. Vathsala Naik"
$ws.Range("A4").Value = "Dr. (Col) M. C. Sharma"
$ws.Range("A5").Value = "Dr. Sanjay Mohanchandra"
